$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet DERMO -> Sheet1
$ws.Name = "Sheet1"

# Header change F1: Pecas -> Pecas.AC
$ws.Range("F1").Value = "Pecas.AC"

# Row 2: becomes numeric values (date text stays as text via quote prefix)
$ws.Range("A2").Value = "'07/07/2023"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 50

# Row 3: becomes plain text strings (no currency/padding)
$ws.Range("A3").Value = "'07/07/2023"
$ws.Range("B3").Value = "'1000.00"
$ws.Range("C3").Value = "'3000.00"
$ws.Range("D3").Value = "'1000.00"
$ws.Range("E3").Value = "'2000.00"
$ws.Range("F3").Value = "'20.0"
$ws.Range("G3").Value = "'1000.00"
$ws.Range("H3").Value = "'66.67"

# Row 4 removed entirely
$ws.Rows.Item(4).Delete()
